$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("A2") "2025/12/10"
Set-TextValue $ws.Range("B2") "5.33"

# Row 8
Set-TextValue $ws.Range("A8") "2025/12/10"
Set-TextValue $ws.Range("B8") "7.89"

# Row 14 (only date changes)
Set-TextValue $ws.Range("A14") "2025/12/10"

# Row 20
Set-TextValue $ws.Range("A20") "2025/12/10"
Set-TextValue $ws.Range("B20") "12.42"

# Row 26
Set-TextValue $ws.Range("A26") "2025/12/10"
Set-TextValue $ws.Range("B26") "10.46"

# Row 32
Set-TextValue $ws.Range("A32") "2025/12/10"
Set-TextValue $ws.Range("B32") "26.32"

# Row 38 (only date changes)
Set-TextValue $ws.Range("A38") "2025/12/10"

# Row 44
Set-TextValue $ws.Range("A44") "2025/12/10"
Set-TextValue $ws.Range("B44") "11.61"

# Row 50
Set-TextValue $ws.Range("A50") "2025/12/10"
Set-TextValue $ws.Range("B50") "12.23"

# Row 56
Set-TextValue $ws.Range("A56") "2025/12/10"
Set-TextValue $ws.Range("B56") "35.84"

# Row 62
Set-TextValue $ws.Range("A62") "2025/12/10"
Set-TextValue $ws.Range("B62") "12.33"

# Row 68
Set-TextValue $ws.Range("A68") "2025/12/10"
Set-TextValue $ws.Range("B68") "14.16"

# Row 74
Set-TextValue $ws.Range("A74") "2025/12/10"
Set-TextValue $ws.Range("B74") "16.84"
